$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2023-12-14 11:01:40", 0.001),
    @("2023-12-14 11:02:11", 0.002),
    @("2023-12-14 11:02:35", 0.0018),
    @("2023-12-14 11:02:53", 0.0002)
)

$startRow = 311
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
